# "search for addclient added and two test customers in excel sheet"
#
# Add two test customers (order no. 4 and 5) to the bottom of the
# "clients_structured" sheet's data table, in the previously-blank row 11
# and the newly-appended row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clients_structured")

# Row 11 (already existed as a blank styled row) and row 12 (new) need the
# same direct formatting used by the rest of the "Name"/"Email"/"Order"
# columns. Seed that formatting first, from cells that already carry it, so
# it sticks once we write the values below.
$ws.Range("B11").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 11 - testcustomer1
$ws.Cells.Item(11, 1).Value = 4
$ws.Cells.Item(11, 2).Value = "testcustomer1"

# Row 12 - testcustomer2
$ws.Cells.Item(12, 1).Value = 5
$ws.Cells.Item(12, 2).Value = "testcustomer2"

# Emails
$ws.Cells.Item(11, 3).Value = "testcustomer1@gmail.com"
$ws.Cells.Item(12, 3).Value = "testcustomer1@gmail.com"

# Orders + quantities
$ws.Cells.Item(11, 4).Value = "Tea"
$ws.Cells.Item(11, 5).Value = 1

$ws.Cells.Item(12, 4).Value = "Water"
$ws.Cells.Item(12, 5).Value = 2

$ws.Range("D17").Select()
